$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("part2")

# Add a new diary entry row (row 4) for exercises 2.5-2.6
$ws.Range("A4").Value = 211103
$ws.Range("B4").Formula = "=30+30"
$ws.Range("C4").Value = "ex 2.5-2.6 and material"

# Move the active selection down to the next empty row, as Excel would
# after data entry in the previous row.
$ws.Range("A5").Select()
